$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 97 (shifts old rows
# 97-162 down to 99-164, matching the target dimension A1:R164).
$ws.Rows("97:98").Insert()

# New row 97: Asterix, 1a (cosecha lavada), $/malla 25 kilos, Región del Maule
$ws.Cells.Item(97, 1).Value = 1
$ws.Cells.Item(97, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(97, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(97, 4).Value = 44978
$ws.Cells.Item(97, 5).Value = 15
$ws.Cells.Item(97, 6).Value = 100114001
$ws.Cells.Item(97, 7).Value = "Papa"
$ws.Cells.Item(97, 8).Value = "Asterix"
$ws.Cells.Item(97, 9).Value = "1a (cosecha lavada)"
$ws.Cells.Item(97, 10).Value = 1000
$ws.Cells.Item(97, 11).Value = 15000
$ws.Cells.Item(97, 12).Value = 16000
$ws.Cells.Item(97, 13).Value = 15500
$ws.Cells.Item(97, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(97, 15).Value = "Región del Maule"
$ws.Cells.Item(97, 16).Value = 620
$ws.Cells.Item(97, 17).Value = 25
$ws.Cells.Item(97, 18).Value = "Hortaliza"

# New row 98: Cardinal, 1a (cosecha), $/malla 25 kilos, Región de Coquimbo
$ws.Cells.Item(98, 1).Value = 1
$ws.Cells.Item(98, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(98, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(98, 4).Value = 44978
$ws.Cells.Item(98, 5).Value = 15
$ws.Cells.Item(98, 6).Value = 100114001
$ws.Cells.Item(98, 7).Value = "Papa"
$ws.Cells.Item(98, 8).Value = "Cardinal"
$ws.Cells.Item(98, 9).Value = "1a (cosecha)"
$ws.Cells.Item(98, 10).Value = 1000
$ws.Cells.Item(98, 11).Value = 14000
$ws.Cells.Item(98, 12).Value = 15000
$ws.Cells.Item(98, 13).Value = 14500
$ws.Cells.Item(98, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(98, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(98, 16).Value = 580
$ws.Cells.Item(98, 17).Value = 25
$ws.Cells.Item(98, 18).Value = "Hortaliza"
